$d = $word.ActiveDocument

# Locate the end of the document body (after the last paragraph, which
# holds the "_GoBack" bookmark) and open up a fresh, empty paragraph
# there to anchor the new content on.
$origCount = $d.Paragraphs.Count
$tail = $d.Paragraphs.Item($origCount).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$insertionPoint = $d.Paragraphs.Item($d.Paragraphs.Count).Range

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = "<w:p $wNs><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>KORIST</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:r><w:t>Od 2003 Python spada u top 10 najpopularnijih  programskih jezika. 2017. je dobio status trečeg najpopularnijeg jezika koji nema sintaksu C-a.</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:r><w:t>Python može služiti kao skripterski jezik za web aplikacijekoristeči web frameworkove kao što je Django.</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:r><w:t>Library poput NumPy, SciPy i Matplotlib omogučuju korištenje pythona u znanstvenoj obradi podataka.</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:r><w:t>Mnogi operativni sustavi dodaju Python kao standardnu komponentu.</w:t></w:r></w:p>" + `
       "<w:p $wNs><w:r><w:t xml:space='preserve'>LibreOffice će uskoro zamjeniti svoj kod u Javi sa Pythonom. </w:t></w:r></w:p>" + `
       "<w:p $wNs/>"

$insertionPoint.InsertXML($xml)

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
